# Update TPM-derived values on the active sheet to reflect the new TPM numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 21.33926
$ws.Range("N2").Value = 64.01778
$ws.Range("O2").Value = 0.4398914187744692
$ws.Range("P2").Value = 0.4398914187744692
$ws.Range("Q2").Value = 222.6590313932666
$ws.Range("R2").Value = 2003.9312825394
$ws.Range("S2").Value = 0.4271826028399997
$ws.Range("T2").Value = 0.4271826028399998

# Row 3
$ws.Range("O3").Value = 0.23906065069302
$ws.Range("P3").Value = 0.23906065069302
$ws.Range("S3").Value = 0.2321539967389683
$ws.Range("T3").Value = 0.2321539967389683

# Row 4
$ws.Range("M4").Value = 10.59425366666667
$ws.Range("N4").Value = 31.782761
$ws.Range("O4").Value = 0.2183918878295978
$ws.Range("P4").Value = 0.2183918878295978
$ws.Range("Q4").Value = 110.5430206930589
$ws.Range("R4").Value = 994.8871862375299
$ws.Range("S4").Value = 0.2120823710135158
$ws.Range("T4").Value = 0.2120823710135158

# Row 5
$ws.Range("M5").Value = 4.979874333333333
$ws.Range("N5").Value = 14.939623
$ws.Range("O5").Value = 0.102656042702913
$ws.Range("P5").Value = 0.102656042702913
$ws.Range("Q5").Value = 51.96122056342111
$ws.Range("R5").Value = 467.6509850707899
$ws.Range("S5").Value = 0.0996902272866745
$ws.Range("T5").Value = 0.09969022728667451

# Row 6
$ws.Range("M6").Value = 21.33926
$ws.Range("N6").Value = 64.01778
$ws.Range("O6").Value = 0.4398914187744692
$ws.Range("P6").Value = 0.4398914187744692
$ws.Range("Q6").Value = 6.62417576772
$ws.Range("R6").Value = 59.61758190947999
$ws.Range("S6").Value = 0.01270881593446946
$ws.Range("T6").Value = 0.01270881593446946

# Row 7
$ws.Range("O7").Value = 0.23906065069302
$ws.Range("P7").Value = 0.23906065069302
$ws.Range("S7").Value = 0.006906653954051677
$ws.Range("T7").Value = 0.006906653954051678

# Row 8
$ws.Range("M8").Value = 10.59425366666667
$ws.Range("N8").Value = 31.782761
$ws.Range("O8").Value = 0.2183918878295978
$ws.Range("P8").Value = 0.2183918878295978
$ws.Range("Q8").Value = 3.288689411714
$ws.Range("R8").Value = 29.598204705426
$ws.Range("S8").Value = 0.006309516816081946
$ws.Range("T8").Value = 0.006309516816081947

# Row 9
$ws.Range("M9").Value = 4.979874333333333
$ws.Range("N9").Value = 14.939623
$ws.Range("O9").Value = 0.102656042702913
$ws.Range("P9").Value = 0.102656042702913
$ws.Range("Q9").Value = 1.545862550302
$ws.Range("R9").Value = 13.912762952718
$ws.Range("S9").Value = 0.002965815416238526
$ws.Range("T9").Value = 0.002965815416238527
